$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet updates ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Blackwater Coal Mine, Australia, M0011, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$row = 2
while ($row -le 23) {
    $cell = $wsData.Range("S$row")
    if ($cell.Text -eq $oldVersion) {
        $cell.Value = $newVersion
    }
    $row = $row + 1
}
